# Auto-generated edit script applying the Hyperion_Profits diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N)
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 53: No Accounting for Waste | Enchanted Electrum Ink
$ws.Range("H53").Value = 4879.773
$ws.Range("I53").Value = 331.08334
$ws.Range("J53").Value = 10338.2
$ws.Range("K53").Value = 331.08334
$ws.Range("L53").Value = 10338.2
$ws.Range("M53").Value = 305.91666
$ws.Range("N53").Value = -11612.2
# row 112: Making Ends Meet | Superior Spiritbond Potion
$ws.Range("H112").Value = 4352.079
$ws.Range("J112").Value = 4352.079
$ws.Range("L112").Value = 13056.237
$ws.Range("N112").Value = -15272.237
# row 113: Amaro Kart | Starch Glue
$ws.Range("H113").Value = 7690.077
$ws.Range("J113").Value = 7108.5557
$ws.Range("L113").Value = 7108.5557
$ws.Range("N113").Value = -13616.5557
# row 135: For Tired Minds | Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 1825.0605
$ws.Range("I135").Value = 713.55
$ws.Range("K135").Value = 6421.95
$ws.Range("M135").Value = -3886.95

$ws = $wb.Worksheets.Item("ARM")
# row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 1408.6
$ws.Range("I61").Value = 1303.1177
$ws.Range("K61").Value = 1303.1177
$ws.Range("M61").Value = -1091.1177
# row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 445515.56
$ws.Range("I122").Value = 1787.3784
$ws.Range("K122").Value = 5362.135200000001
$ws.Range("M122").Value = -2912.135200000001
# row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 1320.2963
$ws.Range("I132").Value = 916.86664
$ws.Range("J132").Value = 3337.4443
$ws.Range("K132").Value = 2750.59992
$ws.Range("L132").Value = 10012.3329
$ws.Range("M132").Value = -220.5999199999997
$ws.Range("N132").Value = -15072.3329
# row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1408.6
$ws.Range("I136").Value = 1303.1177
$ws.Range("K136").Value = 3909.3531
$ws.Range("M136").Value = -1359.3531

$ws = $wb.Worksheets.Item("BSM")
# row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 2440295.8
$ws.Range("I86").Value = 4763066
$ws.Range("J86").Value = 1387.15
$ws.Range("K86").Value = 4763066
$ws.Range("L86").Value = 1387.15
$ws.Range("M86").Value = -4761943
$ws.Range("N86").Value = -3633.15
# row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 2440295.8
$ws.Range("I89").Value = 4763066
$ws.Range("J89").Value = 1387.15
$ws.Range("K89").Value = 23815330
$ws.Range("L89").Value = 6935.75
$ws.Range("M89").Value = -23809714
$ws.Range("N89").Value = -18167.75
# row 94: High Steal | High Steel Nugget
$ws.Range("H94").Value = 4352097
$ws.Range("I94").Value = 5000856.5
$ws.Range("K94").Value = 5000856.5
$ws.Range("M94").Value = -5000405.5
# row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 3073.342
$ws.Range("I134").Value = 934.08
$ws.Range("J134").Value = 7187.3076
$ws.Range("K134").Value = 2802.24
$ws.Range("L134").Value = 21561.9228
$ws.Range("M134").Value = -267.2400000000002
$ws.Range("N134").Value = -26631.9228

$ws = $wb.Worksheets.Item("CRP")
# row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 2335.3333
$ws.Range("I31").Value = 1526.1282
$ws.Range("J31").Value = 4965.25
$ws.Range("K31").Value = 1526.1282
$ws.Range("L31").Value = 4965.25
$ws.Range("M31").Value = -1231.1282
$ws.Range("N31").Value = -5555.25
# row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 2335.3333
$ws.Range("I34").Value = 1526.1282
$ws.Range("J34").Value = 4965.25
$ws.Range("K34").Value = 1526.1282
$ws.Range("L34").Value = 4965.25
$ws.Range("M34").Value = -1324.1282
$ws.Range("N34").Value = -5369.25
# row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 2669.1
$ws.Range("I132").Value = 1834.625
$ws.Range("J132").Value = 6007
$ws.Range("K132").Value = 5503.875
$ws.Range("L132").Value = 18021
$ws.Range("M132").Value = -2973.875
$ws.Range("N132").Value = -23081
# row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 94098.39999999999
$ws.Range("J134").Value = 5998.3335
$ws.Range("L134").Value = 17995.0005
$ws.Range("N134").Value = -23065.0005

$ws = $wb.Worksheets.Item("CUL")
# row 2: Pork Is a Salty Food | Table Salt
$ws.Range("H2").Value = 325.24326
$ws.Range("I2").Value = 138.14285
$ws.Range("J2").Value = 439.13043
$ws.Range("K2").Value = 828.8571000000001
$ws.Range("L2").Value = 2634.78258
$ws.Range("M2").Value = -715.8571000000001
$ws.Range("N2").Value = -2860.78258
# row 34: Fever Pitch | Chamomile Tea
$ws.Range("H34").Value = 1105.875
$ws.Range("I34").Value = 166.33333
$ws.Range("K34").Value = 498.99999
$ws.Range("M34").Value = -414.99999
# row 39: Bloody Good Tart, This | Blood Currant Tart
$ws.Range("H39").Value = 3377.4285
$ws.Range("J39").Value = 3390.3333
$ws.Range("L39").Value = 10170.9999
$ws.Range("N39").Value = -10758.9999
# row 51: The Perks of Life at Sea | Jerked Beef
$ws.Range("H51").Value = 640.2
$ws.Range("I51").Value = 640.2
$ws.Range("K51").Value = 1920.6
$ws.Range("M51").Value = -1460.6
# row 55: Pagan Pastries | Pastry Fish
$ws.Range("H55").Value = 74564.78999999999
$ws.Range("J55").Value = 94817.73
$ws.Range("L55").Value = 284453.19
$ws.Range("N55").Value = -284807.19
# row 56: Culture Club | Crowned Pie
$ws.Range("H56").Value = 19236600
$ws.Range("I56").Value = 19236600
$ws.Range("K56").Value = 19236600
$ws.Range("M56").Value = -19236070
# row 102: Persimmony Snicket | Persimmon Pudding
$ws.Range("H102").Value = 1500
$ws.Range("J102").Value = 1500
$ws.Range("L102").Value = 4500
$ws.Range("N102").Value = -9368

$ws = $wb.Worksheets.Item("GSM")
# row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Range("H80").Value = 1356768.6
$ws.Range("I80").Value = 2711237.5
$ws.Range("J80").Value = 2299.7778
$ws.Range("K80").Value = 2711237.5
$ws.Range("L80").Value = 2299.7778
$ws.Range("M80").Value = -2710239.5
$ws.Range("N80").Value = -4295.7778
# row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Range("H83").Value = 1356768.6
$ws.Range("I83").Value = 2711237.5
$ws.Range("J83").Value = 2299.7778
$ws.Range("K83").Value = 13556187.5
$ws.Range("L83").Value = 11498.889
$ws.Range("M83").Value = -13551195.5
$ws.Range("N83").Value = -21482.889
# row 128: To Fight at Her Side | Manganese Rapier
$ws.Range("H128").Value = 280000
$ws.Range("J128").Value = 280000
$ws.Range("L128").Value = 280000
$ws.Range("N128").Value = -289960
# row 130: Planisphere to Paper | Chondrite Magitek Planisphere
$ws.Range("H130").Value = 60000
$ws.Range("J130").Value = 60000
$ws.Range("L130").Value = 60000
$ws.Range("N130").Value = -70040

$ws = $wb.Worksheets.Item("LTW")
# row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value = 5969.5
$ws.Range("I7").Value = 3616.1667
$ws.Range("K7").Value = 3616.1667
$ws.Range("M7").Value = -3504.1667
# row 36: Campaign in the Membrane | Toadskin Jacket
$ws.Range("H36").Value = 50650
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
# row 46: Supply Side Logic | Boar Leather
$ws.Range("H46").Value = 7009.1
$ws.Range("I46").Value = 6242.6665
$ws.Range("J46").Value = 7636.1816
$ws.Range("K46").Value = 6242.6665
$ws.Range("L46").Value = 7636.1816
$ws.Range("M46").Value = -6054.6665
$ws.Range("N46").Value = -8012.1816
# row 82: Trainin' the Neck | Dragon Leather
$ws.Range("H82").Value = 3269818.8
$ws.Range("I82").Value = 11115000
$ws.Range("J82").Value = 993.25
$ws.Range("K82").Value = 11115000
$ws.Range("L82").Value = 993.25
$ws.Range("M82").Value = -11114639
$ws.Range("N82").Value = -1715.25
# row 85: Training Is Only Skintight (L) | Dragon Leather
$ws.Range("H85").Value = 3269818.8
$ws.Range("I85").Value = 11115000
$ws.Range("J85").Value = 993.25
$ws.Range("K85").Value = 11115000
$ws.Range("L85").Value = 993.25
$ws.Range("M85").Value = -11113752
$ws.Range("N85").Value = -3489.25
# row 93: Hide to Go Seek | Gagana Leather
$ws.Range("H93").Value = 17546738
$ws.Range("I93").Value = 20835952
$ws.Range("J93").Value = 4266
$ws.Range("K93").Value = 20835952
$ws.Range("L93").Value = 4266
$ws.Range("M93").Value = -20834704
$ws.Range("N93").Value = -6762
# row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value = 5969.5
$ws.Range("I126").Value = 3616.1667
$ws.Range("K126").Value = 10848.5001
$ws.Range("M126").Value = -8378.500100000001
# row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 6264.8477
$ws.Range("I132").Value = 6019.3716
$ws.Range("J132").Value = 7045.909
$ws.Range("K132").Value = 18058.1148
$ws.Range("L132").Value = 21137.727
$ws.Range("M132").Value = -15528.1148
$ws.Range("N132").Value = -26197.727
# row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value = 23873.213
$ws.Range("I136").Value = 28117.95
$ws.Range("J136").Value = 3180.125
$ws.Range("K136").Value = 84353.85000000001
$ws.Range("L136").Value = 9540.375
$ws.Range("M136").Value = -81803.85000000001
$ws.Range("N136").Value = -14640.375

$ws = $wb.Worksheets.Item("WVR")
# row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 1675.3043
$ws.Range("I126").Value = 1631.7
$ws.Range("J126").Value = 1966
$ws.Range("K126").Value = 4895.1
$ws.Range("L126").Value = 5898
$ws.Range("M126").Value = -2425.1
$ws.Range("N126").Value = -10838
# row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 62510556
$ws.Range("I132").Value = 76934104
$ws.Range("J132").Value = 8499.666999999999
$ws.Range("K132").Value = 230802312
$ws.Range("L132").Value = 25499.001
$ws.Range("M132").Value = -230799782
$ws.Range("N132").Value = -30559.001
# row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 1861.0682
$ws.Range("I136").Value = 1116.9706
$ws.Range("K136").Value = 3350.9118
$ws.Range("M136").Value = -800.9118000000003

